$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.08840041403804122
$ws.Cells.Item(2, 1).Value = -0.0059999999373196999
$ws.Cells.Item(3, 1).Value = -0.0039999999507749351
$ws.Cells.Item(4, 1).Value = -0.0079999999073017136
$ws.Cells.Item(5, 1).Value = -0.0029999999572760672
$ws.Cells.Item(6, 1).Value = -0.0019999999636208798
$ws.Cells.Item(7, 1).Value = -0.0099999998807502877
$ws.Cells.Item(8, 1).Value = -0.0099999998814923607
$ws.Cells.Item(9, 1).Value = -0.0019999999668072199
$ws.Cells.Item(10, 1).Value = -0.001999999970028199
$ws.Cells.Item(11, 1).Value = -0.00299999995986866
$ws.Cells.Item(12, 1).Value = 0.019696238628236884
$ws.Cells.Item(13, 1).Value = -0.0034999999574969465
$ws.Cells.Item(14, 1).Value = -0.0079999999120294873
$ws.Cells.Item(15, 1).Value = -0.00099999998613320429
$ws.Cells.Item(16, 1).Value = 0.023623603795394743
$ws.Cells.Item(17, 1).Value = -0.0019999999768209875
$ws.Cells.Item(18, 1).Value = -0.0039999999557762678
$ws.Cells.Item(19, 1).Value = -0.0039999999578399503
$ws.Cells.Item(20, 1).Value = 0.0092816714294396974
$ws.Cells.Item(21, 1).Value = -0.0039999999571271871
$ws.Cells.Item(22, 1).Value = -0.0039999999567061906
$ws.Cells.Item(23, 1).Value = -0.0049999999370760051
$ws.Cells.Item(24, 1).Value = -0.019999999774841903
$ws.Cells.Item(25, 1).Value = -0.01999999977139133
$ws.Cells.Item(26, 1).Value = -0.0024999999579815579
$ws.Cells.Item(27, 1).Value = -0.0024999999574473186
$ws.Cells.Item(28, 1).Value = -0.0019999999602404728
$ws.Cells.Item(29, 1).Value = -0.0069999999074861208
$ws.Cells.Item(30, 1).Value = -0.059999999360835332
$ws.Cells.Item(31, 1).Value = 0.054754065990914924
$ws.Cells.Item(32, 1).Value = 0.012282450002768286
$ws.Cells.Item(33, 1).Value = -0.0039999999519508833
